{"js": "// Update the 20x5 table of addition/subtraction equations to the new\n// values from the commit. Each table cell holds exactly one run of text\n// (e.g. \"82-67=15\"); we rewrite the whole table via Table.values so every\n// cell's single text run is updated in place while formatting (fonts,\n// sizes, cell/paragraph properties) is left untouched.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"Expected a table in the document body, found none.\");\n}\n\nconst table = tables.items[0];\ntable.load(\"rowCount,values\");\nawait context.sync();\n\nconst newValues = [\n  [\n    \"37-31=6\",\n    \"51+41=92\",\n    \"58+17=75\",\n    \"14+9=23\",\n    \"13+56=69\"\n  ],\n  [\n    \"76+14=90\",\n    \"84-81=3\",\n    \"26-23=3\",\n    \"20+37=57\",\n    \"70-47=23\"\n  ],\n  [\n    \"62+29=91\",\n    \"40+23=63\",\n    \"10+87=97\",\n    \"97-90=7\",\n    \"52+27=79\"\n  ],\n  [\n    \"66+32=98\",\n    \"43+17=60\",\n    \"73-47=26\",\n    \"11+83=94\",\n    \"48+14=62\"\n  ],\n  [\n    \"24+72=96\",\n    \"89-87=2\",\n    \"50+4=54\",\n    \"43+7=50\",\n    \"95-51=44\"\n  ],\n  [\n    \"69-5=64\",\n    \"37+15=52\",\n    \"93-21=72\",\n    \"92-27=65\",\n    \"85-67=18\"\n  ],\n  [\n    \"77-76=1\",\n    \"13+33=46\",\n    \"17+42=59\",\n    \"4+52=56\",\n    \"28-20=8\"\n  ],\n  [\n    \"86-72=14\",\n    \"84-18=66\",\n    \"29-18=11\",\n    \"76-32=44\",\n    \"78-10=68\"\n  ],\n  [\n    \"11+30=41\",\n    \"69-3=66\",\n    \"9+31=40\",\n    \"57+4=61\",\n    \"26+38=64\"\n  ],\n  [\n    \"16+54=70\",\n    \"40+42=82\",\n    \"43+0=43\",\n    \"53+43=96\",\n    \"70-22=48\"\n  ],\n  [\n    \"21+43=64\",\n    \"76-48=28\",\n    \"6+30=36\",\n    \"50-38=12\",\n    \"97+0=97\"\n  ],\n  [\n    \"15+23=38\",\n    \"83-68=15\",\n    \"46-38=8\",\n    \"13+73=86\",\n    \"77-32=45\"\n  ],\n  [\n    \"97-96=1\",\n    \"68+20=88\",\n    \"76+9=85\",\n    \"5+82=87\",\n    \"16+53=69\"\n  ],\n  [\n    \"99-59=40\",\n    \"77-42=35\",\n    \"17+70=87\",\n    \"83-78=5\",\n    \"87+6=93\"\n  ],\n  [\n    \"98-55=43\",\n    \"25+45=70\",\n    \"59-50=9\",\n    \"90-10=80\",\n    \"21+70=91\"\n  ],\n  [\n    \"91-50=41\",\n    \"35+59=94\",\n    \"74+13=87\",\n    \"14-4=10\",\n    \"33+57=90\"\n  ],\n  [\n    \"27+28=55\",\n    \"91-85=6\",\n    \"36+30=66\",\n    \"3+86=89\",\n    \"96-90=6\"\n  ],\n  [\n    \"79-0=79\",\n    \"30+28=58\",\n    \"61+4=65\",\n    \"78+15=93\",\n    \"88-37=51\"\n  ],\n  [\n    \"18-16=2\",\n    \"59-2=57\",\n    \"29+57=86\",\n    \"32+19=51\",\n    \"2+34=36\"\n  ],\n  [\n    \"30+3=33\",\n    \"27-26=1\",\n    \"17+57=74\",\n    \"9+19=28\",\n    \"98-77=21\"\n  ]\n];\n\nif (table.rowCount !== newValues.length) {\n  throw new Error(`Expected ${newValues.length} rows, found ${table.rowCount}.`);\n}\n\ntable.values = newValues;\nawait context.sync();\n", "ps1": "# Update the 20x5 table of addition/subtraction equations to the new\n# values from the commit. Each table cell holds a single run of text\n# (e.g. \"82-67=15\"); we set each cell's Range.Text directly so only the\n# text changes while cell/paragraph/run formatting is left untouched.\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\n$newValues = @(\n  @(\"37-31=6\",\"51+41=92\",\"58+17=75\",\"14+9=23\",\"13+56=69\"),\n  @(\"76+14=90\",\"84-81=3\",\"26-23=3\",\"20+37=57\",\"70-47=23\"),\n  @(\"62+29=91\",\"40+23=63\",\"10+87=97\",\"97-90=7\",\"52+27=79\"),\n  @(\"66+32=98\",\"43+17=60\",\"73-47=26\",\"11+83=94\",\"48+14=62\"),\n  @(\"24+72=96\",\"89-87=2\",\"50+4=54\",\"43+7=50\",\"95-51=44\"),\n  @(\"69-5=64\",\"37+15=52\",\"93-21=72\",\"92-27=65\",\"85-67=18\"),\n  @(\"77-76=1\",\"13+33=46\",\"17+42=59\",\"4+52=56\",\"28-20=8\"),\n  @(\"86-72=14\",\"84-18=66\",\"29-18=11\",\"76-32=44\",\"78-10=68\"),\n  @(\"11+30=41\",\"69-3=66\",\"9+31=40\",\"57+4=61\",\"26+38=64\"),\n  @(\"16+54=70\",\"40+42=82\",\"43+0=43\",\"53+43=96\",\"70-22=48\"),\n  @(\"21+43=64\",\"76-48=28\",\"6+30=36\",\"50-38=12\",\"97+0=97\"),\n  @(\"15+23=38\",\"83-68=15\",\"46-38=8\",\"13+73=86\",\"77-32=45\"),\n  @(\"97-96=1\",\"68+20=88\",\"76+9=85\",\"5+82=87\",\"16+53=69\"),\n  @(\"99-59=40\",\"77-42=35\",\"17+70=87\",\"83-78=5\",\"87+6=93\"),\n  @(\"98-55=43\",\"25+45=70\",\"59-50=9\",\"90-10=80\",\"21+70=91\"),\n  @(\"91-50=41\",\"35+59=94\",\"74+13=87\",\"14-4=10\",\"33+57=90\"),\n  @(\"27+28=55\",\"91-85=6\",\"36+30=66\",\"3+86=89\",\"96-90=6\"),\n  @(\"79-0=79\",\"30+28=58\",\"61+4=65\",\"78+15=93\",\"88-37=51\"),\n  @(\"18-16=2\",\"59-2=57\",\"29+57=86\",\"32+19=51\",\"2+34=36\"),\n  @(\"30+3=33\",\"27-26=1\",\"17+57=74\",\"9+19=28\",\"98-77=21\")\n)\n\nif ($tbl.Rows.Count -ne $newValues.Count) {\n  throw \"Expected $($newValues.Count) rows, found $($tbl.Rows.Count).\"\n}\n\nfor ($r = 1; $r -le $tbl.Rows.Count; $r++) {\n  $rowValues = $newValues[$r - 1]\n  for ($c = 1; $c -le $rowValues.Count; $c++) {\n    $tbl.Cell($r, $c).Range.Text = $rowValues[$c - 1]\n  }\n}\n\n"}
